# Updates the cryptos price/volume table to the latest scrape snapshot.
# Two coin rows (12/13 and 40/41 and 45/46) also swapped rank order, so
# those rows' Coin/Link/Price/Volume cells are fully rewritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.975.73'
$ws.Range("E2").Value = '  +5.92%  '
$ws.Range("D3").Value = '2.588.66'
$ws.Range("E3").Value = '  +5.76%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.58'
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.55'
$ws.Range("E6").Value = '  +6.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("E8").Value = '  +2.62%  '
$ws.Range("D9").Value = '2.615.06'
$ws.Range("E9").Value = '  +6.79%  '
$ws.Range("E10").Value = '  +5.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.160'
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.32'
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.361'
$ws.Range("E13").Value = '  +3.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.58'
$ws.Range("E14").Value = '  +3.21%  '
$ws.Range("E15").Value = '  +6.12%  '
$ws.Range("D16").Value = '3.060.53'
$ws.Range("E16").Value = '  +5.83%  '
$ws.Range("D17").Value = '65.895.93'
$ws.Range("E17").Value = '  +5.75%  '
$ws.Range("D18").Value = '2.619.02'
$ws.Range("E18").Value = '  +6.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.13'
$ws.Range("E19").Value = '  +4.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.19'
$ws.Range("E20").Value = '  +4.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '355.58'
$ws.Range("E21").Value = '  +10.48%  '
$ws.Range("E22").Value = '  +4.43%  '
$ws.Range("E23").Value = '  +5.38%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.11'
$ws.Range("E25").Value = '  +3.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '66.23'
$ws.Range("E26").Value = '  +1.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '642.07'
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("E28").Value = '  +11.47%  '
$ws.Range("D29").Value = '2.719.55'
$ws.Range("E29").Value = '  +5.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.50'
$ws.Range("E30").Value = '  +6.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.992'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.24'
$ws.Range("E32").Value = '  +5.18%  '
$ws.Range("E33").Value = '  +5.47%  '
$ws.Range("E34").Value = '  +5.95%  '
$ws.Range("E35").Value = '  +8.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.97'
$ws.Range("E37").Value = '  +7.25%  '
$ws.Range("E38").Value = '  +8.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.94'
$ws.Range("E39").Value = '  +9.36%  '
$ws.Range("B40").Value = 'EthereumClassic'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.39'
$ws.Range("E40").Value = '  +5.17%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '156.15'
$ws.Range("E41").Value = '  +3.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.375'
$ws.Range("E42").Value = '  +3.05%  '
$ws.Range("E43").Value = '  +8.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.29'
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₆0316'
$ws.Range("E45").Value = '  +3.17%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '162.32'
$ws.Range("E46").Value = '  +6.82%  '
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.20'
$ws.Range("E48").Value = '  +5.91%  '
$ws.Range("E49").Value = '  +6.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.82'
$ws.Range("E50").Value = '  +9.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.637'
$ws.Range("E51").Value = '  +5.89%  '
